$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 405801.28
$ws.Range("I2").Value = 1212604.4
$ws.Range("J2").Value = 2399.75
$ws.Range("K2").Value = 1212604.4
$ws.Range("L2").Value = 2399.75
$ws.Range("M2").Value = -1212491.4
$ws.Range("N2").Value = -2625.75
$ws.Range("H17").Value = 3572.6667
$ws.Range("J17").Value = 3572.6667
$ws.Range("L17").Value = 10718.0001
$ws.Range("N17").Value = -11054.0001
$ws.Range("H34").Value = 11944.77
$ws.Range("I34").Value = 7040.4287
$ws.Range("K34").Value = 7040.4287
$ws.Range("M34").Value = -6837.4287
$ws.Range("H36").Value = 11944.77
$ws.Range("I36").Value = 7040.4287
$ws.Range("K36").Value = 7040.4287
$ws.Range("M36").Value = -6325.4287
$ws.Range("H43").Value = 4249.9
$ws.Range("J43").Value = 4249.9
$ws.Range("L43").Value = 4249.9
$ws.Range("N43").Value = -4387.9
$ws.Range("H51").Value = 7869.2856
$ws.Range("J51").Value = 12361.667
$ws.Range("L51").Value = 12361.667
$ws.Range("N51").Value = -13329.667
$ws.Range("H53").Value = 986.2727
$ws.Range("I53").Value = 785.7143
$ws.Range("J53").Value = 1337.25
$ws.Range("K53").Value = 785.7143
$ws.Range("L53").Value = 1337.25
$ws.Range("M53").Value = -148.7143
$ws.Range("N53").Value = -2611.25
$ws.Range("H70").Value = 3199
$ws.Range("I70").Value = 2033.3334
$ws.Range("J70").Value = 4073.25
$ws.Range("K70").Value = 6100.0002
$ws.Range("L70").Value = 12219.75
$ws.Range("M70").Value = -5830.0002
$ws.Range("N70").Value = -12759.75
$ws.Range("H73").Value = 3199
$ws.Range("I73").Value = 2033.3334
$ws.Range("J73").Value = 4073.25
$ws.Range("K73").Value = 6100.0002
$ws.Range("L73").Value = 12219.75
$ws.Range("M73").Value = -5164.0002
$ws.Range("N73").Value = -14091.75
$ws.Range("H96").Value = 2602.4375
$ws.Range("I96").Value = 1441.909
$ws.Range("K96").Value = 4325.727000000001
$ws.Range("M96").Value = -2952.727000000001
$ws.Range("H98").Value = 2427.76
$ws.Range("I98").Value = 1850.1818
$ws.Range("J98").Value = 6663.3335
$ws.Range("K98").Value = 1850.1818
$ws.Range("L98").Value = 6663.3335
$ws.Range("M98").Value = -352.1818000000001
$ws.Range("N98").Value = -9659.333500000001
$ws.Range("H106").Value = 4299.875
$ws.Range("I106").Value = 4299.875
$ws.Range("K106").Value = 4299.875
$ws.Range("M106").Value = -3668.875
$ws.Range("H116").Value = 7499.2
$ws.Range("I116").Value = 7499.3335
$ws.Range("J116").Value = 7499
$ws.Range("K116").Value = 7499.3335
$ws.Range("L116").Value = 7499
$ws.Range("M116").Value = -4057.3335
$ws.Range("N116").Value = -14383
$ws.Range("H122").Value = 2427.76
$ws.Range("I122").Value = 1850.1818
$ws.Range("J122").Value = 6663.3335
$ws.Range("K122").Value = 5550.5454
$ws.Range("L122").Value = 19990.0005
$ws.Range("M122").Value = -3100.5454
$ws.Range("N122").Value = -24890.0005
$ws.Range("H132").Value = 4544.154
$ws.Range("I132").Value = 2638.2683
$ws.Range("K132").Value = 7914.804900000001
$ws.Range("M132").Value = -5384.804900000001
$ws.Range("H135").Value = 1987.5714
$ws.Range("I135").Value = 1270.4073
$ws.Range("K135").Value = 11433.6657
$ws.Range("M135").Value = -8898.665700000001
$ws.Range("H138").Value = 4048.6418
$ws.Range("J138").Value = 3913.5107
$ws.Range("L138").Value = 11740.5321
$ws.Range("N138").Value = -22020.5321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 972.5172
$ws.Range("I2").Value = 865
$ws.Range("J2").Value = 1904.3334
$ws.Range("K2").Value = 865
$ws.Range("L2").Value = 1904.3334
$ws.Range("M2").Value = -752
$ws.Range("N2").Value = -2130.3334
$ws.Range("H31").Value = 8606.5
$ws.Range("I31").Value = 7340.5557
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 7340.5557
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -7046.5557
$ws.Range("N31").Value = -20588
$ws.Range("H32").Value = 14124.41
$ws.Range("I32").Value = 9115.896000000001
$ws.Range("J32").Value = 28649.1
$ws.Range("K32").Value = 9115.896000000001
$ws.Range("L32").Value = 28649.1
$ws.Range("M32").Value = -8828.896000000001
$ws.Range("N32").Value = -29223.1
$ws.Range("H45").Value = 2890.9412
$ws.Range("I45").Value = 1219.6
$ws.Range("K45").Value = 1219.6
$ws.Range("M45").Value = -842.5999999999999
$ws.Range("H61").Value = 2780.92
$ws.Range("I61").Value = 1475.1578
$ws.Range("J61").Value = 6915.8335
$ws.Range("K61").Value = 1475.1578
$ws.Range("L61").Value = 6915.8335
$ws.Range("M61").Value = -1263.1578
$ws.Range("N61").Value = -7339.8335
$ws.Range("H74").Value = 2017.8148
$ws.Range("I74").Value = 2018.9231
$ws.Range("J74").Value = 1989
$ws.Range("K74").Value = 2018.9231
$ws.Range("L74").Value = 1989
$ws.Range("M74").Value = -1144.9231
$ws.Range("N74").Value = -3737
$ws.Range("H77").Value = 2017.8148
$ws.Range("I77").Value = 2018.9231
$ws.Range("J77").Value = 1989
$ws.Range("K77").Value = 10094.6155
$ws.Range("L77").Value = 9945
$ws.Range("M77").Value = -5726.6155
$ws.Range("N77").Value = -18681
$ws.Range("H116").Value = 972.5172
$ws.Range("I116").Value = 865
$ws.Range("J116").Value = 1904.3334
$ws.Range("K116").Value = 865
$ws.Range("L116").Value = 1904.3334
$ws.Range("M116").Value = 1429
$ws.Range("N116").Value = -6492.3334
$ws.Range("H132").Value = 1778.8077
$ws.Range("I132").Value = 1701.8723
$ws.Range("J132").Value = 2502
$ws.Range("K132").Value = 5105.6169
$ws.Range("L132").Value = 7506
$ws.Range("M132").Value = -2575.6169
$ws.Range("N132").Value = -12566
$ws.Range("H136").Value = 2780.92
$ws.Range("I136").Value = 1475.1578
$ws.Range("J136").Value = 6915.8335
$ws.Range("K136").Value = 4425.4734
$ws.Range("L136").Value = 20747.5005
$ws.Range("M136").Value = -1875.4734
$ws.Range("N136").Value = -25847.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 972.5172
$ws.Range("I3").Value = 865
$ws.Range("J3").Value = 1904.3334
$ws.Range("K3").Value = 865
$ws.Range("L3").Value = 1904.3334
$ws.Range("M3").Value = -751
$ws.Range("N3").Value = -2132.3334
$ws.Range("H20").Value = 1301.76
$ws.Range("I20").Value = 1176.8182
$ws.Range("J20").Value = 2218
$ws.Range("K20").Value = 1176.8182
$ws.Range("L20").Value = 2218
$ws.Range("M20").Value = -929.8181999999999
$ws.Range("N20").Value = -2712
$ws.Range("H99").Value = 2319.1333
$ws.Range("I99").Value = 1752.8462
$ws.Range("K99").Value = 1752.8462
$ws.Range("M99").Value = -254.8462
$ws.Range("H102").Value = 11837.167
$ws.Range("I102").Value = 11837.167
$ws.Range("K102").Value = 11837.167
$ws.Range("M102").Value = -8592.166999999999
$ws.Range("H105").Value = 2062.8462
$ws.Range("I105").Value = 1628.75
$ws.Range("J105").Value = 2757.4
$ws.Range("K105").Value = 1628.75
$ws.Range("L105").Value = 2757.4
$ws.Range("M105").Value = 118.25
$ws.Range("N105").Value = -6251.4
$ws.Range("H110").Value = 44351
$ws.Range("J110").Value = 44351
$ws.Range("L110").Value = 44351
$ws.Range("N110").Value = -52531
$ws.Range("H132").Value = 112496.75
$ws.Range("J132").Value = 112496.75
$ws.Range("L132").Value = 112496.75
$ws.Range("N132").Value = -122616.75
$ws.Range("H133").Value = 100649.7
$ws.Range("I133").Value = 50000
$ws.Range("J133").Value = 103315.48
$ws.Range("K133").Value = 50000
$ws.Range("L133").Value = 103315.48
$ws.Range("N133").Value = -113435.48
$ws.Range("H134").Value = 1777.283
$ws.Range("I134").Value = 1813.6666
$ws.Range("J134").Value = 849.5
$ws.Range("K134").Value = 5440.9998
$ws.Range("L134").Value = 2548.5
$ws.Range("M134").Value = -2905.9998
$ws.Range("N134").Value = -7618.5
$ws.Range("M133").Value = -44940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1213.7273
$ws.Range("I16").Value = 1021.1429
$ws.Range("J16").Value = 1550.75
$ws.Range("K16").Value = 1021.1429
$ws.Range("L16").Value = 1550.75
$ws.Range("M16").Value = -734.1429000000001
$ws.Range("N16").Value = -2124.75
$ws.Range("H31").Value = 2696.244
$ws.Range("I31").Value = 1704.1471
$ws.Range("J31").Value = 7515
$ws.Range("K31").Value = 1704.1471
$ws.Range("L31").Value = 7515
$ws.Range("M31").Value = -1409.1471
$ws.Range("N31").Value = -8105
$ws.Range("H34").Value = 2696.244
$ws.Range("I34").Value = 1704.1471
$ws.Range("J34").Value = 7515
$ws.Range("K34").Value = 1704.1471
$ws.Range("L34").Value = 7515
$ws.Range("M34").Value = -1502.1471
$ws.Range("N34").Value = -7919
$ws.Range("H58").Value = 2489.5
$ws.Range("I58").Value = 2245.5
$ws.Range("K58").Value = 2245.5
$ws.Range("M58").Value = -2042.5
$ws.Range("H59").Value = 95500
$ws.Range("J59").Value = 100000
$ws.Range("L59").Value = 100000
$ws.Range("N59").Value = -102290
$ws.Range("H86").Value = 21668.709
$ws.Range("J86").Value = 15002.333
$ws.Range("L86").Value = 15002.333
$ws.Range("N86").Value = -17248.333
$ws.Range("H89").Value = 21668.709
$ws.Range("J89").Value = 15002.333
$ws.Range("L89").Value = 75011.66500000001
$ws.Range("N89").Value = -86243.66500000001
$ws.Range("H99").Value = 3336.6924
$ws.Range("J99").Value = 4566.3335
$ws.Range("L99").Value = 4566.3335
$ws.Range("N99").Value = -7562.3335
$ws.Range("H113").Value = 1213.7273
$ws.Range("I113").Value = 1021.1429
$ws.Range("J113").Value = 1550.75
$ws.Range("K113").Value = 1021.1429
$ws.Range("L113").Value = 1550.75
$ws.Range("M113").Value = 1148.8571
$ws.Range("N113").Value = -5890.75
$ws.Range("H122").Value = 962.35486
$ws.Range("I122").Value = 995.7778
$ws.Range("K122").Value = 2987.3334
$ws.Range("M122").Value = -537.3334
$ws.Range("H126").Value = 3336.6924
$ws.Range("J126").Value = 4566.3335
$ws.Range("L126").Value = 13699.0005
$ws.Range("N126").Value = -18639.0005
$ws.Range("H132").Value = 3132.56
$ws.Range("I132").Value = 3024.3635
$ws.Range("K132").Value = 9073.0905
$ws.Range("M132").Value = -6543.0905
$ws.Range("H134").Value = 2207.7878
$ws.Range("I134").Value = 2174.4517
$ws.Range("K134").Value = 6523.355100000001
$ws.Range("M134").Value = -3988.355100000001
$ws.Range("H136").Value = 2489.5
$ws.Range("I136").Value = 2245.5
$ws.Range("K136").Value = 6736.5
$ws.Range("M136").Value = -4186.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2980.7693
$ws.Range("J48").Value = 2980.7693
$ws.Range("L48").Value = 8942.3079
$ws.Range("N48").Value = -9442.3079
$ws.Range("H74").Value = 956.5
$ws.Range("I74").Value = 956.5
$ws.Range("K74").Value = 2869.5
$ws.Range("H77").Value = 956.5
$ws.Range("I77").Value = 956.5
$ws.Range("K77").Value = 8608.5
$ws.Range("M74").Value = -1808.5
$ws.Range("M77").Value = -3304.5
$ws.Range("H92").Value = 632.6896400000001
$ws.Range("I92").Value = 525.6875
$ws.Range("J92").Value = 764.38464
$ws.Range("K92").Value = 1577.0625
$ws.Range("L92").Value = 2293.15392
$ws.Range("M92").Value = -329.0625
$ws.Range("N92").Value = -4789.15392
$ws.Range("H107").Value = 436.0435
$ws.Range("I107").Value = 412.8
$ws.Range("J107").Value = 442.5
$ws.Range("K107").Value = 1238.4
$ws.Range("L107").Value = 1327.5
$ws.Range("M107").Value = 681.5999999999999
$ws.Range("N107").Value = -5167.5
$ws.Range("H122").Value = 730.0714
$ws.Range("I122").Value = 614.7778
$ws.Range("J122").Value = 937.6
$ws.Range("K122").Value = 5533.000199999999
$ws.Range("L122").Value = 8438.4
$ws.Range("M122").Value = -3083.000199999999
$ws.Range("N122").Value = -13338.4
$ws.Range("H123").Value = 4665.6665
$ws.Range("I123").Value = 2599
$ws.Range("J123").Value = 14999
$ws.Range("K123").Value = 7797
$ws.Range("L123").Value = 44997
$ws.Range("M123").Value = -5347
$ws.Range("N123").Value = -49897
$ws.Range("H131").Value = 4360111
$ws.Range("I131").Value = 3643.3333
$ws.Range("J131").Value = 5013581.5
$ws.Range("K131").Value = 10929.9999
$ws.Range("L131").Value = 15040744.5
$ws.Range("M131").Value = -5889.999899999999
$ws.Range("N131").Value = -15050824.5
$ws.Range("H134").Value = 2527.5
$ws.Range("I134").Value = 2527.5
$ws.Range("K134").Value = 7582.5
$ws.Range("M134").Value = -2512.5
$ws.Range("H139").Value = 1534
$ws.Range("I139").Value = 1039.6666
$ws.Range("K139").Value = 3118.9998
$ws.Range("M139").Value = 2021.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9557
$ws.Range("I43").Value = 1039.0834
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 1039.0834
$ws.Range("L43").Value = 30000
$ws.Range("H57").Value = 11513.375
$ws.Range("I57").Value = 11513.375
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 11513.375
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -10693.375
$ws.Range("M43").Value = -888.0834
$ws.Range("N43").Value = -30302
$ws.Range("H63").Value = 37444
$ws.Range("J63").Value = 37444
$ws.Range("L63").Value = 37444
$ws.Range("N63").Value = -38816
$ws.Range("H66").Value = 37444
$ws.Range("J66").Value = 37444
$ws.Range("L66").Value = 112332
$ws.Range("N66").Value = -119196
$ws.Range("H80").Value = 2026.4706
$ws.Range("I80").Value = 1695.6364
$ws.Range("J80").Value = 2633
$ws.Range("K80").Value = 1695.6364
$ws.Range("L80").Value = 2633
$ws.Range("M80").Value = -697.6364000000001
$ws.Range("N80").Value = -4629
$ws.Range("H83").Value = 2026.4706
$ws.Range("I83").Value = 1695.6364
$ws.Range("J83").Value = 2633
$ws.Range("K83").Value = 8478.182000000001
$ws.Range("L83").Value = 13165
$ws.Range("M83").Value = -3486.182000000001
$ws.Range("N83").Value = -23149
$ws.Range("H97").Value = 891.7838
$ws.Range("J97").Value = 1018.6667
$ws.Range("L97").Value = 1018.6667
$ws.Range("N97").Value = -2010.6667
$ws.Range("H113").Value = 143997.42
$ws.Range("I113").Value = 84413.75
$ws.Range("K113").Value = 84413.75
$ws.Range("M113").Value = -82243.75
$ws.Range("H132").Value = 2222.2126
$ws.Range("I132").Value = 1303.8649
$ws.Range("J132").Value = 5620.1
$ws.Range("K132").Value = 3911.5947
$ws.Range("L132").Value = 16860.3
$ws.Range("M132").Value = -1381.5947
$ws.Range("N132").Value = -21920.3
$ws.Range("H139").Value = 74010
$ws.Range("J139").Value = 74010
$ws.Range("L139").Value = 74010
$ws.Range("N139").Value = -84290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7335.4707
$ws.Range("I7").Value = 7533.533
$ws.Range("J7").Value = 5850
$ws.Range("K7").Value = 7533.533
$ws.Range("L7").Value = 5850
$ws.Range("M7").Value = -7421.533
$ws.Range("N7").Value = -6074
$ws.Range("H16").Value = 772.55
$ws.Range("I16").Value = 646.8
$ws.Range("K16").Value = 646.8
$ws.Range("M16").Value = -476.8
$ws.Range("H42").Value = 33249
$ws.Range("I42").Value = 27999
$ws.Range("K42").Value = 27999
$ws.Range("M42").Value = -27436
$ws.Range("H46").Value = 10103.454
$ws.Range("J46").Value = 4998
$ws.Range("L46").Value = 4998
$ws.Range("N46").Value = -5374
$ws.Range("H49").Value = 33249
$ws.Range("I49").Value = 27999
$ws.Range("K49").Value = 27999
$ws.Range("M49").Value = -27852
$ws.Range("H55").Value = 3233.9092
$ws.Range("I55").Value = 3228.7144
$ws.Range("K55").Value = 3228.7144
$ws.Range("M55").Value = -3055.7144
$ws.Range("H82").Value = 2018.5405
$ws.Range("I82").Value = 1039.4375
$ws.Range("J82").Value = 2764.524
$ws.Range("K82").Value = 1039.4375
$ws.Range("L82").Value = 2764.524
$ws.Range("M82").Value = -678.4375
$ws.Range("N82").Value = -3486.524
$ws.Range("H85").Value = 2018.5405
$ws.Range("I85").Value = 1039.4375
$ws.Range("J85").Value = 2764.524
$ws.Range("K85").Value = 1039.4375
$ws.Range("L85").Value = 2764.524
$ws.Range("M85").Value = 208.5625
$ws.Range("N85").Value = -5260.523999999999
$ws.Range("H94").Value = 55000
$ws.Range("J94").Value = 55000
$ws.Range("L94").Value = 55000
$ws.Range("H122").Value = 4496.2085
$ws.Range("I122").Value = 3499.875
$ws.Range("J122").Value = 4994.375
$ws.Range("K122").Value = 10499.625
$ws.Range("L122").Value = 14983.125
$ws.Range("M122").Value = -8049.625
$ws.Range("N122").Value = -19883.125
$ws.Range("H126").Value = 7335.4707
$ws.Range("I126").Value = 7533.533
$ws.Range("J126").Value = 5850
$ws.Range("K126").Value = 22600.599
$ws.Range("L126").Value = 17550
$ws.Range("M126").Value = -20130.599
$ws.Range("N126").Value = -22490
$ws.Range("H132").Value = 4099.864
$ws.Range("I132").Value = 3932
$ws.Range("J132").Value = 5163
$ws.Range("K132").Value = 11796
$ws.Range("L132").Value = 15489
$ws.Range("M132").Value = -9266
$ws.Range("N132").Value = -20549
$ws.Range("H136").Value = 3603
$ws.Range("I136").Value = 3635.125
$ws.Range("J136").Value = 3474.5
$ws.Range("K136").Value = 10905.375
$ws.Range("L136").Value = 10423.5
$ws.Range("M136").Value = -8355.375
$ws.Range("N136").Value = -15523.5
$ws.Range("N94").Value = -56352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 56100
$ws.Range("J93").Value = 56100
$ws.Range("L93").Value = 56100
$ws.Range("N93").Value = -61092
$ws.Range("H107").Value = 968.6667
$ws.Range("I107").Value = 503.83334
$ws.Range("J107").Value = 1898.3334
$ws.Range("K107").Value = 1511.50002
$ws.Range("L107").Value = 5695.0002
$ws.Range("M107").Value = 408.4999800000001
$ws.Range("N107").Value = -9535.0002
$ws.Range("H132").Value = 1966.2941
$ws.Range("I132").Value = 1966.2941
$ws.Range("K132").Value = 5898.8823
$ws.Range("M132").Value = -3368.8823
$ws.Range("H136").Value = 4789.375
$ws.Range("I136").Value = 4822.5
$ws.Range("J136").Value = 4716.5
$ws.Range("K136").Value = 14467.5
$ws.Range("L136").Value = 14149.5
$ws.Range("M136").Value = -11917.5
$ws.Range("N136").Value = -19249.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N57").ClearContents()
